$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# AMAZONAS (row 2)
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 1

# APURIMAC (row 4)
$ws.Range("B4").Value = 3

# CAJAMARCA (row 7)
$ws.Range("C7").Value = 6

# CALLAO (row 8)
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0

# JUNIN (row 13)
$ws.Range("E13").Value = 1

# LAMBAYEQUE (row 15)
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1

# MOQUEGUA (row 20)
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1

# TACNA (row 25)
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 1

# UCAYALI (row 27)
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
